# Refresh the cryptocurrency price/volume figures on the active sheet.
# D column = "Price", E column = "Volume(1h)" (percentage change text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '37.847.58'
$ws.Range("E2").Value = '  +6.33%  '
$ws.Range("D3").Value = "'" + '2.056.60'
$ws.Range("E3").Value = '  +3.73%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = "'" + '252.93'
$ws.Range("E5").Value = '  +4.33%  '
$ws.Range("E6").Value = '  +2.31%  '
$ws.Range("D7").Value = "'" + '65.32'
$ws.Range("E7").Value = '  +14.10%  '
$ws.Range("D8").Value = "'" + '1.00'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = "'" + '60.95'
$ws.Range("E9").Value = '  +2.10%  '
$ws.Range("D10").Value = "'" + '0.379'
$ws.Range("E10").Value = '  +5.33%  '
$ws.Range("D11").Value = "'" + '0.0765'
$ws.Range("E11").Value = '  +4.81%  '
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").Value = "'" + '0.922'
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("D14").Value = "'" + '15.21'
$ws.Range("E14").Value = '  +8.66%  '
$ws.Range("D15").Value = "'" + '2.349.93'
$ws.Range("E15").Value = '  +3.25%  '
$ws.Range("D16").Value = "'" + '20.78'
$ws.Range("E16").Value = '  +20.91%  '
$ws.Range("D17").Value = "'" + '5.56'
$ws.Range("E17").Value = '  +6.12%  '
$ws.Range("D18").Value = "'" + '2.031.40'
$ws.Range("E18").Value = '  +2.49%  '
$ws.Range("D19").Value = "'" + '37.695.89'
$ws.Range("E19").Value = '  +6.09%  '
$ws.Range("D20").Value = "'" + '74.13'
$ws.Range("E20").Value = '  +4.99%  '
$ws.Range("D21").Value = "'" + '0.0₃0881'
$ws.Range("E21").Value = '  +5.19%  '
$ws.Range("D22").Value = "'" + '5.36'
$ws.Range("E22").Value = '  +5.97%  '
$ws.Range("D23").Value = "'" + '239.36'
$ws.Range("E23").Value = '  +2.56%  '
$ws.Range("D24").Value = "'" + '2.70'
$ws.Range("E24").Value = '  +15.70%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  +5.25%  '
$ws.Range("D27").Value = "'" + '9.67'
$ws.Range("E27").Value = '  +5.92%  '
$ws.Range("D28").Value = "'" + '160.12'
$ws.Range("E28").Value = '  -2.06%  '
$ws.Range("D29").Value = "'" + '20.01'
$ws.Range("E29").Value = '  +2.98%  '
$ws.Range("D30").Value = "'" + '0.115'
$ws.Range("E30").Value = '  +28.71%  '
$ws.Range("E31").Value = '  +2.84%  '
$ws.Range("E32").Value = '  +8.99%  '
$ws.Range("D33").Value = "'" + '1.21'
$ws.Range("E33").Value = '  +6.94%  '
$ws.Range("E34").Value = '  +11.13%  '
$ws.Range("E35").Value = '  +5.08%  '
$ws.Range("D36").Value = "'" + '2.45'
$ws.Range("E36").Value = '  +3.58%  '
$ws.Range("E37").Value = '  +3.75%  '
$ws.Range("E38").Value = '  -0.13%  '
$ws.Range("D39").Value = "'" + '6.10'
$ws.Range("E39").Value = '  +24.04%  '
$ws.Range("E40").Value = '  +17.51%  '
$ws.Range("D41").Value = "'" + '2.81'
$ws.Range("E41").Value = '  +24.79%  '
$ws.Range("E42").Value = '  +4.28%  '
$ws.Range("E43").Value = '  +4.71%  '
$ws.Range("E44").Value = '  +3.57%  '
$ws.Range("E45").Value = '  +5.95%  '
$ws.Range("D46").Value = "'" + '17.05'
$ws.Range("E46").Value = '  +10.63%  '
$ws.Range("E47").Value = '  +7.70%  '
$ws.Range("D48").Value = "'" + '95.51'
$ws.Range("E48").Value = '  +5.01%  '
$ws.Range("D49").Value = "'" + '1.415.21'
$ws.Range("E49").Value = '  +2.81%  '
$ws.Range("E50").Value = '  +2.46%  '
$ws.Range("D51").Value = "'" + '47.36'
$ws.Range("E51").Value = '  +3.51%  '
